# Update the "Started" (Yes/No) flag in column C of the "by Coach" sheet
# for the coaches/players whose status changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

$updates = @{
    5  = "Yes"
    6  = "Yes"
    9  = "No"
    10 = "No"
    17 = "Yes"
    22 = "No"
    43 = "Yes"
    45 = "No"
    51 = "Yes"
    53 = "No"
    54 = "No"
    56 = "Yes"
    76 = "No"
    78 = "Yes"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

# Reset the frozen-pane scroll position / selection back to the top of the sheet.
$ws.Range("A2").Select()
